# Fix minor bugs in configuration writer
# Updates the values of the transformer loading-percent result block (B2:H25)
# with corrected values produced by the fixed configuration writer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 27.67192888141118
$ws.Cells.Item(2, 3).Value = 23.49754641113499
$ws.Cells.Item(2, 4).Value = 4.428107757948536
$ws.Cells.Item(2, 5).Value = 3.345900308208168
$ws.Cells.Item(2, 6).Value = 3.384954215364001
$ws.Cells.Item(2, 7).Value = 1.304850381317482
$ws.Cells.Item(2, 8).Value = 2.342917425459456
$ws.Cells.Item(3, 2).Value = 26.52339124900794
$ws.Cells.Item(3, 3).Value = 22.59165154910132
$ws.Cells.Item(3, 4).Value = 2.342930740777862
$ws.Cells.Item(3, 5).Value = 2.318844236779696
$ws.Cells.Item(3, 6).Value = 2.34197080088003
$ws.Cells.Item(3, 7).Value = 0.2638765009624271
$ws.Cells.Item(3, 8).Value = 2.340652271129138
$ws.Cells.Item(4, 2).Value = 25.47187834739725
$ws.Cells.Item(4, 3).Value = 21.65654622576773
$ws.Cells.Item(4, 4).Value = 2.341130862470552
$ws.Cells.Item(4, 5).Value = 1.292797567134669
$ws.Cells.Item(4, 6).Value = 1.302165880704331
$ws.Cells.Item(4, 7).Value = 2.342138922042389
$ws.Cells.Item(4, 8).Value = 2.338924569364795
$ws.Cells.Item(5, 2).Value = 25.36206422289756
$ws.Cells.Item(5, 3).Value = 21.70239389623367
$ws.Cells.Item(5, 4).Value = 1.302360395417469
$ws.Cells.Item(5, 5).Value = 3.344001214354973
$ws.Cells.Item(5, 6).Value = 1.301964356597287
$ws.Cells.Item(5, 7).Value = 1.302781259094242
$ws.Cells.Item(5, 8).Value = 0.2646766362126541
$ws.Cells.Item(6, 2).Value = 24.77451109222771
$ws.Cells.Item(6, 3).Value = 21.0424170190098
$ws.Cells.Item(6, 4).Value = 2.33985028450725
$ws.Cells.Item(6, 5).Value = 0.2677523358254872
$ws.Cells.Item(6, 6).Value = 2.338950626655075
$ws.Cells.Item(6, 7).Value = 1.302440255972398
$ws.Cells.Item(6, 8).Value = 2.337735040658941
$ws.Cells.Item(7, 2).Value = 25.49764281472686
$ws.Cells.Item(7, 3).Value = 21.54889398327296
$ws.Cells.Item(7, 4).Value = 5.459959003146573
$ws.Cells.Item(7, 5).Value = 2.318145781413409
$ws.Cells.Item(7, 6).Value = 4.41820460380544
$ws.Cells.Item(7, 7).Value = 2.342723416892716
$ws.Cells.Item(7, 8).Value = 3.377328055448976
$ws.Cells.Item(8, 2).Value = 26.74674662495746
$ws.Cells.Item(8, 3).Value = 22.67284913424233
$ws.Cells.Item(8, 4).Value = 3.384132350752711
$ws.Cells.Item(8, 5).Value = 3.345026048059967
$ws.Cells.Item(8, 6).Value = 3.382581542273305
$ws.Cells.Item(8, 7).Value = 1.304160103138903
$ws.Cells.Item(8, 8).Value = 5.459265959033391
$ws.Cells.Item(9, 2).Value = 30.35713475138244
$ws.Cells.Item(9, 3).Value = 25.770130725514
$ws.Cells.Item(9, 4).Value = 3.394888338634906
$ws.Cells.Item(9, 5).Value = 6.43223858013834
$ws.Cells.Item(9, 6).Value = 1.306081399738332
$ws.Cells.Item(9, 7).Value = 7.579741332876252
$ws.Cells.Item(9, 8).Value = 2.347808827026818
$ws.Cells.Item(10, 2).Value = 43.14783282689503
$ws.Cells.Item(10, 3).Value = 36.19245974954005
$ws.Cells.Item(10, 4).Value = 3.434748034845658
$ws.Cells.Item(10, 5).Value = 4.391577769511869
$ws.Cells.Item(10, 6).Value = 7.665747152655189
$ws.Cells.Item(10, 7).Value = 5.558227957003028
$ws.Cells.Item(10, 8).Value = 10.83060369914447
$ws.Cells.Item(11, 2).Value = 54.07041101900032
$ws.Cells.Item(11, 3).Value = 45.58135484474686
$ws.Cells.Item(11, 4).Value = 5.612435321586476
$ws.Cells.Item(11, 5).Value = 11.66570456161316
$ws.Cells.Item(11, 6).Value = 7.747856681136822
$ws.Cells.Item(11, 7).Value = 4.545853928745167
$ws.Cells.Item(11, 8).Value = 5.597561579235927
$ws.Cells.Item(12, 2).Value = 72.62995236996102
$ws.Cells.Item(12, 3).Value = 60.48622802877691
$ws.Cells.Item(12, 4).Value = 16.70734712004668
$ws.Cells.Item(12, 5).Value = 12.78103684560337
$ws.Cells.Item(12, 6).Value = 9.002524034679595
$ws.Cells.Item(12, 7).Value = 7.935244832618292
$ws.Cells.Item(12, 8).Value = 13.35570364242531
$ws.Cells.Item(13, 2).Value = 87.15118921846668
$ws.Cells.Item(13, 3).Value = 72.14410008889118
$ws.Cells.Item(13, 4).Value = 11.41386590317339
$ws.Cells.Item(13, 5).Value = 14.94636299714753
$ws.Cells.Item(13, 6).Value = 18.09017932972794
$ws.Cells.Item(13, 7).Value = 18.17289641190914
$ws.Cells.Item(13, 8).Value = 13.58235030369825
$ws.Cells.Item(14, 2).Value = 88.51991021341182
$ws.Cells.Item(14, 3).Value = 73.0795212844795
$ws.Cells.Item(14, 4).Value = 14.7971741918863
$ws.Cells.Item(14, 5).Value = 9.698104379710291
$ws.Cells.Item(14, 6).Value = 18.12185444395588
$ws.Cells.Item(14, 7).Value = 15.95600040811754
$ws.Cells.Item(14, 8).Value = 14.72053351423004
$ws.Cells.Item(15, 2).Value = 87.30390055764565
$ws.Cells.Item(15, 3).Value = 72.36873693892068
$ws.Cells.Item(15, 4).Value = 13.6517841424199
$ws.Cells.Item(15, 5).Value = 15.99942727412624
$ws.Cells.Item(15, 6).Value = 13.62092815289322
$ws.Cells.Item(15, 7).Value = 15.92726942867577
$ws.Cells.Item(15, 8).Value = 14.69593428855574
$ws.Cells.Item(16, 2).Value = 87.18210141189483
$ws.Cells.Item(16, 3).Value = 72.39518271779393
$ws.Cells.Item(16, 4).Value = 8.057515578752398
$ws.Cells.Item(16, 5).Value = 18.1039162667512
$ws.Cells.Item(16, 6).Value = 14.73336469003992
$ws.Cells.Item(16, 7).Value = 17.04413321883194
$ws.Cells.Item(16, 8).Value = 14.69157298312121
$ws.Cells.Item(17, 2).Value = 87.55203157959885
$ws.Cells.Item(17, 3).Value = 72.2555097165835
$ws.Cells.Item(17, 4).Value = 18.15040130971568
$ws.Cells.Item(17, 5).Value = 17.05053909485994
$ws.Cells.Item(17, 6).Value = 18.10928346992362
$ws.Cells.Item(17, 7).Value = 13.69696729746186
$ws.Cells.Item(17, 8).Value = 23.64020315383295
$ws.Cells.Item(18, 2).Value = 74.02395712917837
$ws.Cells.Item(18, 3).Value = 61.14185062833186
$ws.Cells.Item(18, 4).Value = 23.38241023284229
$ws.Cells.Item(18, 5).Value = 18.01626806607335
$ws.Cells.Item(18, 6).Value = 16.72563791372906
$ws.Cells.Item(18, 7).Value = 15.68533579557764
$ws.Cells.Item(18, 8).Value = 19.98164845797066
$ws.Cells.Item(19, 2).Value = 50.3846366282033
$ws.Cells.Item(19, 3).Value = 41.28237962668724
$ws.Cells.Item(19, 4).Value = 18.47944346243567
$ws.Cells.Item(19, 5).Value = 17.87082482159128
$ws.Cells.Item(19, 6).Value = 21.67542530910544
$ws.Cells.Item(19, 7).Value = 26.04488424600707
$ws.Cells.Item(19, 8).Value = 21.64063160187417
$ws.Cells.Item(20, 2).Value = 43.77102320366958
$ws.Cells.Item(20, 3).Value = 36.14672940143365
$ws.Cells.Item(20, 4).Value = 15.14138254225008
$ws.Cells.Item(20, 5).Value = 21.98312696080114
$ws.Cells.Item(20, 6).Value = 15.12462474986573
$ws.Cells.Item(20, 7).Value = 19.42637042301817
$ws.Cells.Item(20, 8).Value = 24.67840459508772
$ws.Cells.Item(21, 2).Value = 37.13992716945414
$ws.Cells.Item(21, 3).Value = 30.63017148416026
$ws.Cells.Item(21, 4).Value = 22.45278371960438
$ws.Cells.Item(21, 5).Value = 26.08105560667831
$ws.Cells.Item(21, 6).Value = 12.91304296474789
$ws.Cells.Item(21, 7).Value = 19.29352325790306
$ws.Cells.Item(21, 8).Value = 22.40463300696685
$ws.Cells.Item(22, 2).Value = 34.12379181051589
$ws.Cells.Item(22, 3).Value = 28.49988243895644
$ws.Cells.Item(22, 4).Value = 18.14325227883835
$ws.Cells.Item(22, 5).Value = 23.98860564550115
$ws.Cells.Item(22, 6).Value = 11.81236356215059
$ws.Cells.Item(22, 7).Value = 14.99460593718256
$ws.Cells.Item(22, 8).Value = 10.75042245098614
$ws.Cells.Item(23, 2).Value = 33.07044996865772
$ws.Cells.Item(23, 3).Value = 27.5902708961242
$ws.Cells.Item(23, 4).Value = 14.96010283660451
$ws.Cells.Item(23, 5).Value = 15.71219077430775
$ws.Cells.Item(23, 6).Value = 8.647212231564032
$ws.Cells.Item(23, 7).Value = 11.81527519925887
$ws.Cells.Item(23, 8).Value = 11.78524865362885
$ws.Cells.Item(24, 2).Value = 30.59867842236875
$ws.Cells.Item(24, 3).Value = 25.7195387341162
$ws.Cells.Item(24, 4).Value = 7.580768965266447
$ws.Cells.Item(24, 5).Value = 10.54845716292746
$ws.Cells.Item(24, 6).Value = 8.622064108235044
$ws.Cells.Item(24, 7).Value = 3.399356249578878
$ws.Cells.Item(24, 8).Value = 11.75112772520713
$ws.Cells.Item(25, 2).Value = 29.30012398247544
$ws.Cells.Item(25, 3).Value = 24.81025050823691
$ws.Cells.Item(25, 4).Value = 3.391898384034924
$ws.Cells.Item(25, 5).Value = 5.402167924135912
$ws.Cells.Item(25, 6).Value = 4.432925601835993
$ws.Cells.Item(25, 7).Value = 3.393779541215449
$ws.Cells.Item(25, 8).Value = 5.471883964783395
